# buildings-template.xlsx: "update the example" fix.
#
# Row 2 becomes a highlighted (yellow) example row showing a sample
# building name ("May An Phu") next to an "Example" label in column B.
# Row 3 becomes the (red-highlighted) first real input row, left blank
# for the user to fill in. Column B is now part of the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: example data -------------------------------------------------
# (Column B's "Example" label is registered as a shared string before
# column A's sample value, matching the order the sheet was authored in.)
$ws.Range("B2").Value = "Example"
$ws.Range("A2").Value = "May An Phu"

# Highlight the whole example row yellow.
$ws.Range("A2:B2").Interior.Color = 65535   # RGB(255,255,0)

# --- Row 3: first real (blank) entry row, highlighted red ---------------
$ws.Range("A3").Interior.Color = 255        # RGB(255,0,0)

# Put the selection on the first cell the user should actually fill in.
$ws.Range("A3").Select() | Out-Null
